$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (2022-02-14 / serial 44971) is inserted as
# the new most-recent row. It becomes the new row 4, pushing every
# subsequent row (old row 4 .. old row 35) down by one, so the sheet grows
# from 35 to 36 data-bearing rows (A1:R35 -> A1:R36).
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44971
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100112044
$ws.Range("G4").Value = "Perejil"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 2500
$ws.Range("L4").Value = 2800
$ws.Range("M4").Value = 2671
$ws.Range("N4").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O4").Value = "Región de Arica y Parinacota"
$ws.Range("P4").Value = 1336
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = "Hortaliza"
